# Rename sheets in sample excel data files
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("XlsxDataSource1").Name = "Data Sheet 1"
$wb.Worksheets.Item("XlsxDataSource2").Name = "Data Sheet 2"
$wb.Worksheets.Item("XlsxDataSource3").Name = "Data Sheet 3"
